$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.170203924179077
$ws.Range("B1").Value = 2.22244930267334
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.33198070526123
$ws.Range("E1").Value = 1.228568315505981
